$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.263.61'
$ws.Range("E2").Value = '  +1.07%  '
$ws.Range("D3").Value = '2.643.85'
$ws.Range("E3").Value = '  +2.86%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.56'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.586'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.30%  '
$ws.Range("D9").Value = '2.644.19'
$ws.Range("E9").Value = '  +2.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.107'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.68'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.152'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.77%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.354'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.42'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.49%  '
$ws.Range("D15").Value = '3.114.32'
$ws.Range("E15").Value = '  +2.60%  '
$ws.Range("D16").Value = '63.145.86'
$ws.Range("E16").Value = '  +0.98%  '
$ws.Range("E17").Value = '  +0.47%  '
$ws.Range("D18").Value = '2.636.39'
$ws.Range("E18").Value = '  +2.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.38'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '339.08'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.37'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.02%  '
$ws.Range("E22").Value = '  +1.04%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.48%  '
$ws.Range("E25").Value = '  +5.64%  '
$ws.Range("E26").Value = '  +2.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.165'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.41%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("E29").Value = '  +2.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.82'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '523.43'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +15.10%  '
$ws.Range("E32").Value = '  +13.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.98'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.91%  '
$ws.Range("D34").Value = '0.0₃0807'
$ws.Range("E34").Value = '  +0.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '174.46'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.93'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +11.65%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.94%  '
$ws.Range("E40").Value = '  +7.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '171.61'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.06'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.18%  '
$ws.Range("E44").Value = '  +1.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0559'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.630'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0961'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.38%  '
$ws.Range("E49").Value = '  +2.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.81%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.71'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.48%  '
